# Updated cryptos list on Mon May 15 13:08:47 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds plain-text numeric-looking values (e.g. "0.8800",
# "27.725.42") that must stay exactly as typed rather than being
# re-interpreted/normalized as numbers. Mark the whole column as Text
# before writing so the literal strings round-trip unchanged.
$ws.Range("D2:D51").NumberFormat = "@"

function Set-Price($row, $price) {
    $ws.Cells.Item($row, 4).Value = $price
}

function Set-Volume($row, $vol) {
    $ws.Cells.Item($row, 5).Value = "  $vol  "
}

# Row 2 - Bitcoin
Set-Price 2 "27.725.42"
Set-Volume 2 "-0.07%"

# Row 3 - Ethereum
Set-Price 3 "1.848.45"
Set-Volume 3 "-0.89%"

# Row 4 - TetherUSD
Set-Volume 4 "-2.16%"

# Row 5 - BNB
Set-Price 5 "320.08"
Set-Volume 5 "-1.03%"

# Row 6 - USDC
Set-Price 6 "1.011"
Set-Volume 6 "-2.10%"

# Row 7 - XRP
Set-Price 7 "0.4301"
Set-Volume 7 "-2.69%"

# Row 8 - Cardano
Set-Price 8 "0.3739"
Set-Volume 8 "-1.53%"

# Row 9 - Dogecoin
Set-Price 9 "0.07363"
Set-Volume 9 "-1.49%"

# Row 10 - Polygon
Set-Price 10 "0.8800"

# Row 11 - Solana
Set-Price 11 "21.66"
Set-Volume 11 "-0.50%"

# Row 12 - WrappedEther
Set-Price 12 "1.862.79"
Set-Volume 12 "-0.47%"

# Row 13 - Chainlink
Set-Price 13 "6.728"
Set-Volume 13 "-0.48%"

# Row 14 - Polkadot
Set-Price 14 "5.455"
Set-Volume 14 "-1.73%"

# Row 15 - TRON
Set-Price 15 "0.07152"
Set-Volume 15 "-0.98%"

# Row 16 - Litecoin
Set-Price 16 "87.93"
Set-Volume 16 "+4.36%"

# Row 17 - BinanceUSD
Set-Price 17 "1.015"
Set-Volume 17 "-2.37%"

# Row 18 - ShibaInu
Set-Price 18 "0.000008999"
Set-Volume 18 "-1.42%"

# Row 19 - Dai
Set-Price 19 "1.012"
Set-Volume 19 "-2.11%"

# Row 20 - Avalanche
Set-Price 20 "15.43"
Set-Volume 20 "-0.88%"

# Row 21 - WrappedBTC
Set-Price 21 "27.732.06"
Set-Volume 21 "-0.10%"

# Row 22 - Uniswap
Set-Price 22 "5.238"
Set-Volume 22 "-1.43%"

# Row 23 - Cosmos
Set-Volume 23 "-1.56%"

# Row 24 - WrappedliquidstakedEther2.0
Set-Price 24 "2.079.75"
Set-Volume 24 "-1.02%"

# Row 25 - Toncoin
Set-Price 25 "2.005"
Set-Volume 25 "-0.33%"

# Row 26 - Monero
Set-Price 26 "155.73"
Set-Volume 26 "-1.82%"

# Row 27 - EthereumClassic
Set-Price 27 "18.62"
Set-Volume 27 "-1.15%"

# Row 28 - LidoDAOToken
Set-Price 28 "2.138"
Set-Volume 28 "+7.23%"

# Row 29 - InternetComputer(DFINITY)
Set-Price 29 "5.386"
Set-Volume 29 "+1.29%"

# Row 30 - BitcoinCash
Set-Price 30 "119.72"
Set-Volume 30 "+1.48%"

# Row 31 - Stellar
Set-Price 31 "0.08955"
Set-Volume 31 "-1.11%"

# Row 32 - ARBITRUM
Set-Price 32 "1.238"
Set-Volume 32 "+1.66%"

# Row 33 - ImmutableX
Set-Price 33 "0.7794"
Set-Volume 33 "+0.01%"

# Row 34 - Filecoin
Set-Price 34 "4.569"
Set-Volume 34 "-0.23%"

# Row 35 - HuobiToken
Set-Price 35 "2.911"
Set-Volume 35 "-3.58%"

# Row 36 - Frax
Set-Price 36 "1.012"
Set-Volume 36 "-2.26%"

# Row 37 - TrustWalletToken
Set-Volume 37 "-0.92%"

# Row 38 - Hedera
Set-Price 38 "0.05343"
Set-Volume 38 "+0.04%"

# Row 39 - VeChain
Set-Price 39 "0.01972"
Set-Volume 39 "-0.85%"

# Row 40 - FraxShare
Set-Price 40 "7.280"
Set-Volume 40 "+6.02%"

# Row 41 - MXToken
Set-Price 41 "2.882"
Set-Volume 41 "+0.26%"

# Row 42 - TheSandbox
Set-Price 42 "0.5155"
Set-Volume 42 "-1.00%"

# Row 43 - Algorand
Set-Volume 43 "-0.39%"

# Row 44 - Aptos
Set-Price 44 "8.809"
Set-Volume 44 "+1.74%"

# Row 45 - Quant
Set-Price 45 "109.45"
Set-Volume 45 "-0.82%"

# Row 46 - EnergySwap
Set-Price 46 "10.69"
Set-Volume 46 "+0.29%"

# Row 47 - Decentraland
Set-Price 47 "0.4744"
Set-Volume 47 "+0.57%"

# Rows 48 & 49 swap: NEARProtocol and Cronos swap ranking positions.
$ws.Cells.Item(48, 2).Value = "Cronos"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-Price 48 "0.06489"
Set-Volume 48 "-1.88%"

$ws.Cells.Item(49, 2).Value = "NEARProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-Price 49 "1.699"
Set-Volume 49 "-1.13%"

# Row 50 - PaxDollar
Set-Price 50 "1.012"
Set-Volume 50 "-2.32%"

# Row 51 - RenderToken
Set-Price 51 "1.866"
Set-Volume 51 "-3.18%"
